$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Intro sentence: insert "Hello YouTubers, here is the guy with the Swiss
#    accent. " before "If you plan to start with " and wrap a fresh (empty)
#    "_GoBack" bookmark right at the boundary between the two sentences.
# ---------------------------------------------------------------------------
$introTarget = $d.Content
$introTarget.Find.Execute("If you plan to start with", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Hello YouTubers, here is the guy with the Swiss accent. If you plan to start with", 2) | Out-Null

$introBoundary = $d.Content
$introBoundary.Find.Execute("If you plan to start with", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$introBookmarkRange = $d.Range($introBoundary.Start, $introBoundary.Start)
$d.Bookmarks.Add("_GoBack", $introBookmarkRange) | Out-Null

# ---------------------------------------------------------------------------
# 2) Collapse the three runs ("...an installed " / "Raspbian" / " jessie
#    release...") - including the spell-check proofErr wrappers around
#    "Raspbian" - into one contiguous run of plain text.
# ---------------------------------------------------------------------------
$raspbianOld = "After assembling the hat and the Pi, we have to create an SD card with the needed software. We start with an SD card with an installed Raspbian jessie release. You find a link in the description on how to proceed."
$raspbianTarget = $d.Content
$raspbianTarget.Find.Execute($raspbianOld, $false, $false, $false, $false, $false, $true, 1, $false, $raspbianOld, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Drop the stale <w:lastRenderedPageBreak/> in front of "So, let's play
#    around with it" without merging that run into its neighbours. A
#    throw-away bookmark dropped right after the run's text acts as a run
#    boundary "pin" so the paragraph's later runs stay untouched; the
#    re-write of the run's own text (identity replace) flushes the stale
#    page-break marker. The scratch bookmark is then removed again.
# ---------------------------------------------------------------------------
$playSentence = [string][char]0x2019
$playSentence = "So, let" + $playSentence + "s play around with it"

$playFind = $d.Content
$playFind.Find.Execute($playSentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pinRange = $d.Range($playFind.End, $playFind.End)
$d.Bookmarks.Add("zzzScratchPin", $pinRange) | Out-Null

$playRewrite = $d.Content
$playRewrite.Find.Execute($playSentence, $false, $false, $false, $false, $false, $true, 1, $false, $playSentence, 2) | Out-Null

$d.Bookmarks("zzzScratchPin").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove the three trailing "link dump" paragraphs (the Raspberry Pi
#    imaging link, the hackster.io hyperlink and the AliExpress link) plus
#    the bookmark that used to wrap them, leaving the "I hope, this video…"
#    paragraph as the last paragraph of the body.
# ---------------------------------------------------------------------------
$closing = $d.Content
$closing.Find.Execute("I hope, this video was useful or at least interesting for you. If true, then like. Bye", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteStart = $closing.End + 1
$deleteEnd = $d.Content.End
$tail = $d.Range($deleteStart, $deleteEnd)
$tail.Delete() | Out-Null
